# Append a new job listing to the top of the data rows (row 7) on the
# "ランサーズ" sheet, push the previous row 7 down to row 8, and refresh
# the "取得日時" (fetched-at) timestamp for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-17 06:36:16"

# --- 1) Copy the current row 7 (the oldest listing) down to row 8 -------
$ws.Range("A7:H7").Copy($ws.Range("A8:H8"))

# --- 2) Overwrite row 7 with the newly scraped listing -------------------
$ws.Range("D7").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5473394"

# --- 3) Refresh the "取得日時" timestamp on every data row (2-8) ---------
$ws.Range("A2").Value = $newTimestamp
$ws.Range("A3").Value = $newTimestamp
$ws.Range("A4").Value = $newTimestamp
$ws.Range("A5").Value = $newTimestamp
$ws.Range("A6").Value = $newTimestamp
$ws.Range("A7").Value = $newTimestamp
$ws.Range("A8").Value = $newTimestamp

# --- 4) Rebuild the hyperlinks for column F so the rIds line up with the
#        new layout (F7 -> new URL, F8 -> the URL that moved down) -------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5473383")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5217096")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5473147")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5473146")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5473234")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5473394")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5473181")
